$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 2052
$ws.Range("F7").Value = 7823
$ws.Range("F8").Value = 252
$ws.Range("F9").Value = 33
$ws.Range("F12").Value = 1752
$ws.Range("F13").Value = 1520
$ws.Range("F16").Value = 3878
$ws.Range("F17").Value = 5951
$ws.Range("F18").Value = 676
$ws.Range("F21").Value = 1219
$ws.Range("F22").Value = 407
$ws.Range("F23").Value = 6141
$ws.Range("F24").Value = 344
$ws.Range("F25").Value = 52
$ws.Range("F26").Value = 4174
$ws.Range("F27").Value = 697
$ws.Range("F28").Value = 1922
$ws.Range("F29").Value = 1153
$ws.Range("F31").Value = 13
$ws.Range("F32").Value = 7
$ws.Range("F35").Value = 35
$ws.Range("F37").Value = 1143
$ws.Range("F40").Value = 93
$ws.Range("F41").Value = 397
$ws.Range("F42").Value = 147
$ws.Range("F43").Value = 1116
$ws.Range("F45").Value = 60
$ws.Range("F46").Value = 30

# Sheet: 演出
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 4
$ws.Range("F11").Value = 665
$ws.Range("F12").Value = 356
$ws.Range("F20").Value = 162
$ws.Range("F21").Value = 147
$ws.Range("F22").Value = 67
$ws.Range("F25").Value = 89

# Sheet: 本地生活
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F8").Value = 3077
$ws.Range("F9").Value = 895
$ws.Range("F10").Value = 1045
$ws.Range("F11").Value = 1214
$ws.Range("F12").Value = 1544

# Sheet: 全部类型
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 3077
$ws.Range("F8").Value = 2052
$ws.Range("F9").Value = 7823
$ws.Range("F10").Value = 33
$ws.Range("F11").Value = 895
$ws.Range("F15").Value = 1752
$ws.Range("F16").Value = 1520
$ws.Range("F17").Value = 1214
$ws.Range("F19").Value = 665
$ws.Range("F21").Value = 1544
$ws.Range("F22").Value = 3878
$ws.Range("F23").Value = 356
$ws.Range("F25").Value = 676
$ws.Range("F28").Value = 1219
$ws.Range("F29").Value = 408
$ws.Range("F30").Value = 6142
$ws.Range("F31").Value = 344
$ws.Range("F32").Value = 697
$ws.Range("F33").Value = 1922
$ws.Range("F34").Value = 1153
$ws.Range("F37").Value = 162
$ws.Range("F39").Value = 67
$ws.Range("F42").Value = 93
$ws.Range("F43").Value = 397
$ws.Range("F44").Value = 1116
